$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 304
$ws.Range("F1").Value = 2546
$ws.Range("G1").Value = 5779

$ws.Range("E2").Value = 300
$ws.Range("F2").Value = 4544
$ws.Range("G2").Value = 17489

$ws.Range("H13").Select()
